$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.555.53"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "2.303.37"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.59"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.00"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.62"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.43"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.964"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.34"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "2.646.63"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "2.298.91"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "42.482.65"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.41"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000106"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.34"
$ws.Range("E21").Value = "  +28.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.31"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.58"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "273.77"
$ws.Range("E24").Value = "  +5.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("E25").Value = "  -4.62%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.82"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.33"
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.16"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  +3.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0875"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").Value = "  -8.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.59"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0360"
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.78"
$ws.Range("E39").Value = "  +4.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.52"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.47"
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "94.16"
$ws.Range("E43").Value = "  -5.31%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.225"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.08"
$ws.Range("E46").Value = "  -3.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.55"
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "79.90"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.95"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.27"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("D51").Value = "1.605.92"
$ws.Range("E51").Value = "  +4.17%  "
